$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New pl_mw.xlsx line-flow results (380 kV case) for rows 2-25 (A2:A25 = 0..23).
# Columns D, L, N are always 0 and are left untouched.
$data = @(
    @(0.2548275340367354,0.05100019544823908,0.5837379712630479,2.014861140361404,0.4002651927567982,0.5773205202246245,0.4603573165992358,0.04075061024878224,0.2669875471203795,0.4116067872690437,1.901172913236792),
    @(0.2224060481534025,0.04731724755653488,0.5733558073701417,2.007384071406719,0.405329823849307,0.5830368421796521,0.4669261400233662,0.04078715264090071,0.2335551542479664,0.3878202468631144,1.923667207578859),
    @(0.2024267842359109,0.0450433318645338,0.5672541858230957,2.003925152746646,0.4087383533255107,0.5867937445209463,0.4712230754862468,0.04081817130700216,0.2129334820532591,0.3733402334961653,1.938622448692144),
    @(0.1942674612005533,0.04411359653747127,0.5648364898840015,2.002800331980865,0.4102024090473435,0.5883868636693883,0.4730403769000677,0.04083297782120354,0.2045069160809447,0.3674712665587236,1.945004271801061),
    @(0.192911565283822,0.04395902960483511,0.5644391894450393,2.002630754284098,0.410450044720065,0.5886551544301781,0.4733461397415368,0.04083556750565087,0.2031063151932813,0.3664986564677903,1.946081325285569),
    @(0.2023168152869914,0.04503080557962846,0.5672213013749854,2.003908830138059,0.4087577943286504,0.5868149782013532,0.4712473159852575,0.04081836221096857,0.2128199310626684,0.3732609535531068,1.938707352539694),
    @(0.2436639575179811,0.04973296179205988,0.5801015804401999,2.012048120590805,0.4019494025532637,0.5792402473777756,0.4625675075200597,0.04076143389398013,0.2554799162597021,0.4033794239959789,1.908691497613248),
    @(0.3241499366449148,0.05885171298893965,0.6075240464848264,2.036992546899782,0.3909732725698802,0.566345259464903,0.4476396313478777,0.04071754021184759,0.338368125273746,0.4634227939562621,1.858911068301879),
    @(0.3828963921730235,0.06548635687910576,0.628990383837575,2.060802677564226,0.3843626676071921,0.5580637616733455,0.4379496403777665,0.04072614477448333,0.398774590049527,0.5081249533812908,1.827881549641702),
    @(0.4095331179052266,0.06849002154322648,0.6390425070479608,2.07282705502206,0.3816721693352818,0.5545548381997705,0.4338191179042781,0.04073883782801246,0.4261437471649288,0.528587068439144,1.814971105925721),
    @(0.4196066837079968,0.06962529392416172,0.6428901979744523,2.077551973268243,0.3806989892391073,0.5532632302441982,0.4322949366205044,0.04074489876068554,0.4364913990845309,0.536353525599651,1.810255709237978),
    @(0.4174377573821459,0.06938088944335163,0.642059699389506,2.076526748116692,0.3809065487661982,0.5535397492136624,0.432621418342892,0.04074353777038375,0.4342635888844768,0.5346800895794956,1.811263533892173),
    @(0.4103621440336269,0.06858346457831033,0.6393582343444848,2.073212338877624,0.3815911896239896,0.5544478324372974,0.4336929214710104,0.04073931137145692,0.4269953880722142,0.5292256628215455,1.814579688755984),
    @(0.4060263910517108,0.0680947367338689,0.6377088676682803,2.071204507780479,0.382016500633874,0.5550088967928843,0.4343544533601591,0.04073688569563316,0.4225412437941145,0.5258869895339444,1.816633531903079),
    @(0.3811538086025337,0.06528976322373126,0.6283392203872467,2.06004086230547,0.3845448804038796,0.5582982764539111,0.4382251683665483,0.04072549121805658,0.3969836829905944,0.5067902295212932,1.82874954523102),
    @(0.3658724710269325,0.06356524745493175,0.6226646782938872,2.053497885075799,0.3861771538218335,0.5603823734813957,0.440670839802241,0.04072074380956003,0.3812763130590326,0.4951072357231254,1.836491126920208),
    @(0.3570748670322246,0.06257199393517965,0.6194278493738352,2.049846831815273,0.3871457968659087,0.5616054105887045,0.4421036420390916,0.04071884000513926,0.3722315286544813,0.4883994498220687,1.84105729062702),
    @(0.3540947631933875,0.06223546431041882,0.61833655797723,2.048629934799564,0.3874788779420086,0.5620236871131539,0.4425932482023818,0.04071833765526556,0.3691673667001965,0.4861303740158789,1.842622786729223),
    @(0.3675000469290524,0.06374896612199166,0.6232659474693207,2.05418277420253,0.3860003102621974,0.5601580008662239,0.4404077905342909,0.04072116366868705,0.3829494621191714,0.4963496747305385,1.835655283120417),
    @(0.412440785979868,0.068817746501324,0.6401506035447539,2.07418120554037,0.3813888541018358,0.5541800983367509,0.4333771101571404,0.04074051878859564,0.4291306863909483,0.5308272771863756,1.813600942938962),
    @(0.4417350208381379,0.07211792214171453,0.6514256276200712,2.0882512414845,0.3786411543796362,0.5504897004930669,0.4290150842145657,0.04076047607947686,0.4592165899967426,0.5534645474661062,1.80019859868392),
    @(0.4261074005335956,0.07035772897349091,0.6453860140267125,2.080650310335955,0.3800832634538338,0.5524395275457223,0.4313218488278618,0.0407491584964319,0.4431681862416497,0.5413731938086102,1.807259060989153),
    @(0.3667642576918126,0.06366591253467391,0.6229940341098796,2.05387279134186,0.3860801670879965,0.5602593622811298,0.4405266317781322,0.04072097127839491,0.3821930762612737,0.4957879402358998,1.836032808576618),
    @(0.3024425506540638,0.05639605883851573,0.599873878132648,2.029281967782481,0.3936877512670094,0.5696241574497662,0.4514537761857724,0.04072220838811447,0.3160293033705841,0.4470754509051318,1.871404985216202)
)

$cols = @("B","C","E","F","G","H","I","J","K","M","O")

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $row).Value = $data[$i][$j]
    }
}
